$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-12 Wednesday", "2025-11-13 Thursday"),
    @("50×21=1050", "65×87=5655"),
    @("24×70=1680", "12×77=924"),
    @("93×31=2883", "60×58=3480"),
    @("89×80=7120", "59×96=5664"),
    @("68×67=4556", "80×18=1440"),
    @("58×58=3364", "42×29=1218"),
    @("45×80=3600", "58×15=870"),
    @("23×21=483", "22×49=1078"),
    @("55×36=1980", "78×58=4524"),
    @("59×40=2360", "77×20=1540"),
    @("71×15=1065", "93×81=7533"),
    @("66×39=2574", "93×16=1488"),
    @("98×64=6272", "56×95=5320"),
    @("31×81=2511", "98×62=6076"),
    @("18×26=468", "90×88=7920"),
    @("97×68=6596", "91×17=1547"),
    @("95×63=5985", "92×76=6992"),
    @("63×32=2016", "84×55=4620"),
    @("73×17=1241", "91×25=2275"),
    @("86×28=2408", "64×24=1536"),
    @("45×39=1755", "15×27=405"),
    @("97×79=7663", "42×20=840"),
    @("31×51=1581", "53×53=2809"),
    @("90×64=5760", "51×17=867"),
    @("34×73=2482", "39×41=1599")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
